# Revised Prepare.supply to handle cst hubs and terminals
#
# The "C1" label used for the two schedule rows is replaced with "CST"
# (the sheet now distinguishes CST hubs/terminals instead of a generic "C1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "CST"
$ws.Range("A4").Value = "CST"

# Move the active selection (in the frozen bottom-left pane) to C5.
$ws.Range("C5").Select()
